$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (style) used by the other header cells (B1:H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data values for columns I (I0) and J (IF), rows 2-20
$data = @(
    @(8, 9),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(6, 7),
    @(4, 6),
    @(7, 8),
    @(6, 6),
    @(7, 8),
    @(6, 7),
    @(6, 6),
    @(6, 7),
    @(9, 9),
    @(9, 9),
    @(4, 5),
    @(8, 8),
    @(8, 8),
    @(4, 4),
    @(4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
